$d = $word.ActiveDocument
$replacements = @(
    @("2025-02-17 Monday", "2025-02-18 Tuesday"),
    @("9+76=85", "20+79=99"),
    @("38+9=47", "31+20=51"),
    @("3+76=79", "60+4=64"),
    @("75-47=28", "93-13=80"),
    @("27+32=59", "72-17=55"),
    @("90-17=73", "24-23=1"),
    @("11+4=15", "59-15=44"),
    @("36+15=51", "90-35=55"),
    @("16+10=26", "25-8=17"),
    @("18+53=71", "96-4=92"),
    @("84-79=5", "42+54=96"),
    @("87-1=86", "41-27=14"),
    @("45+3=48", "90+8=98"),
    @("21-3=18", "70-23=47"),
    @("83-79=4", "22+4=26"),
    @("82-34=48", "8+66=74"),
    @("14+15=29", "72-23=49"),
    @("90-65=25", "73-27=46"),
    @("59-0=59", "91+1=92"),
    @("20+50=70", "93-50=43"),
    @("20-13=7", "1+69=70"),
    @("8+54=62", "48+7=55"),
    @("7+14=21", "7+50=57"),
    @("86-76=10", "51+14=65"),
    @("98-87=11", "46+10=56"),
    @("22-15=7", "18-2=16"),
    @("77-60=17", "97-81=16"),
    @("45+14=59", "49-6=43"),
    @("5+17=22", "24+1=25"),
    @("95-4=91", "76+18=94"),
    @("37-23=14", "55-1=54"),
    @("98-57=41", "28+63=91"),
    @("67+9=76", "24+15=39"),
    @("90-80=10", "76+1=77"),
    @("11+12=23", "10+20=30"),
    @("35+1=36", "74-48=26"),
    @("31+65=96", "24+7=31"),
    @("7+63=70", "18+70=88"),
    @("28-21=7", "29+35=64"),
    @("1+5=6", "12+70=82"),
    @("7+0=7", "46+8=54"),
    @("17+0=17", "68-36=32"),
    @("66+3=69", "21-9=12"),
    @("96-62=34", "51-49=2"),
    @("32+7=39", "29-5=24"),
    @("28-15=13", "41-18=23"),
    @("58-27=31", "67+5=72"),
    @("79-30=49", "65-21=44"),
    @("2+62=64", "59+12=71"),
    @("61+5=66", "91-53=38"),
    @("19+37=56", "63-56=7"),
    @("93-23=70", "66-3=63"),
    @("39+52=91", "71+16=87"),
    @("4+34=38", "57-31=26"),
    @("74-58=16", "16-4=12"),
    @("68-41=27", "38+59=97"),
    @("92-29=63", "10+40=50"),
    @("18+54=72", "98-36=62"),
    @("56+42=98", "59+29=88"),
    @("53+21=74", "53+8=61"),
    @("91-66=25", "52-5=47"),
    @("31+21=52", "6+59=65"),
    @("13+43=56", "5-3=2"),
    @("84-35=49", "20+17=37"),
    @("25+7=32", "99-71=28"),
    @("47-43=4", "63+36=99"),
    @("80+11=91", "74-49=25"),
    @("67+31=98", "98-2=96"),
    @("84-70=14", "52+32=84"),
    @("21-15=6", "36+54=90"),
    @("93-92=1", "46-29=17"),
    @("89-38=51", "3+96=99"),
    @("92-84=8", "52-1=51"),
    @("96-12=84", "46+5=51"),
    @("4+62=66", "66-57=9"),
    @("40+8=48", "35-13=22"),
    @("29-9=20", "78-22=56"),
    @("1+46=47", "63-20=43"),
    @("33+44=77", "99-65=34"),
    @("46+32=78", "38-20=18"),
    @("18-4=14", "55+22=77"),
    @("46-22=24", "98-21=77"),
    @("51-38=13", "52-42=10"),
    @("34+58=92", "34-23=11"),
    @("89-30=59", "10+16=26"),
    @("84-59=25", "56-6=50"),
    @("10-3=7", "34-1=33"),
    @("45+18=63", "48-36=12"),
    @("90-60=30", "40+57=97"),
    @("97-19=78", "33-30=3"),
    @("16-15=1", "77+17=94"),
    @("42+44=86", "43-39=4"),
    @("66+2=68", "53-9=44"),
    @("11+80=91", "90-19=71"),
    @("16+8=24", "48+15=63"),
    @("13+0=13", "37+42=79"),
    @("9+64=73", "44+40=84"),
    @("54-25=29", "73-31=42"),
    @("90-69=21", "62-19=43"),
    @("40+6=46", "81-7=74"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}
Write-Host "Done"